# Commit: "adding full run for ZEV Jan R2-4 and modifying files for consistency in R2-4"
#
# Semantic change:
#  - Column E (s1cDNAPreparer) for rows 2-41: value "J.Plaggenberg" -> "J.PLAGGENBERG"
#    and the font color switches from theme-color font to an explicit RGB-black font
#    (this is the pre-existing style index 4 in the workbook).
#  - Column B (rnaPreparer) for rows 22-41: value "J.Plaggenberg" -> "J.PLAGGENBERG"
#    (style/format unchanged).
#  - Active cell selection on the sheet moves from M14 to L19.
#  - The Excel window position shifts (cosmetic / not data - best effort only).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newName = "J.PLAGGENBERG"

# Update column E (s1cDNAPreparer) for rows 2 through 41: new text + font color
# (re-applies the existing "black RGB" font format already present in the workbook
# as style index 4, rather than the theme-color font used before).
for ($r = 2; $r -le 41; $r++) {
    $cell = $ws.Cells.Item($r, 5)   # column E
    $cell.Value = $newName
    $cell.Font.Color = 0
}

# Update column B (rnaPreparer) for rows 22 through 41
for ($r = 22; $r -le 41; $r++) {
    $cell = $ws.Cells.Item($r, 2)   # column B
    $cell.Value = $newName
}

# Update the active selection to L19 (was M14)
$ws.Range("L19").Select()

# Best-effort: move the application window (cosmetic workbook view state)
$win = $wb.Windows.Item(1)
$win.Left = 16080
$win.Top = 460
